# Apply latest cryptos snapshot values (prices + 1h volume deltas) to the sheet.
# Some Price values are plain decimals (e.g. "600.44") which Excel COM would
# otherwise auto-convert to a Double; a leading apostrophe (quote-prefix) keeps
# them as literal text, matching the source data exactly (no float rounding).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.598.55'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '2.671.35'
$ws.Range("E3").Value = '  -1.03%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''600.44'
$ws.Range("E5").Value = '  -1.41%  '

$ws.Range("D6").Value = '''156.92'
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '''0.623'
$ws.Range("E8").Value = '  +5.85%  '

$ws.Range("E9").Value = '  +4.16%  '

$ws.Range("E10").Value = '  -0.52%  '

$ws.Range("E11").Value = '  -2.75%  '

$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("E13").Value = '  -3.13%  '

$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").Value = '3.149.57'
$ws.Range("E15").Value = '  -1.31%  '

$ws.Range("D16").Value = '65.484.17'
$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").Value = '2.657.92'
$ws.Range("E17").Value = '  -1.20%  '

$ws.Range("D18").Value = '''12.78'
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").Value = '''4.80'
$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").Value = '''7.57'
$ws.Range("E20").Value = '  -0.58%  '

$ws.Range("D21").Value = '''351.14'
$ws.Range("E21").Value = '  -2.37%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").Value = '''69.55'
$ws.Range("E23").Value = '  -0.98%  '

$ws.Range("D24").Value = '''0.0000111'
$ws.Range("E24").Value = '  +5.30%  '

$ws.Range("D25").Value = '''9.67'
$ws.Range("E25").Value = '  -0.82%  '

$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = '''1.60'
$ws.Range("E27").Value = '  -5.66%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.166'
$ws.Range("E28").Value = '  -1.68%  '

$ws.Range("D29").Value = '''8.10'
$ws.Range("E29").Value = '  -1.37%  '

$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '''534.56'
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.14'
$ws.Range("E32").Value = '  -2.81%  '

$ws.Range("E33").Value = '  -2.17%  '

$ws.Range("E34").Value = '  -2.65%  '

$ws.Range("E35").Value = '  +1.00%  '

$ws.Range("E36").Value = '  -2.00%  '

$ws.Range("D37").Value = '''20.47'
$ws.Range("E37").Value = '  -1.40%  '

$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("D39").Value = '''158.10'
$ws.Range("E39").Value = '  -3.24%  '

$ws.Range("E40").Value = '  -2.83%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").Value = '''164.59'
$ws.Range("E42").Value = '  -3.02%  '

$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("D44").Value = '''2.34'
$ws.Range("E44").Value = '  +2.11%  '

$ws.Range("D45").Value = '''0.0609'
$ws.Range("E45").Value = '  -0.61%  '

$ws.Range("E46").Value = '  -2.86%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0259'
$ws.Range("E47").Value = '  -2.61%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.640'
$ws.Range("E48").Value = '  -2.74%  '

$ws.Range("D49").Value = '0.0₆0261'
$ws.Range("E49").Value = '  +15.02%  '

$ws.Range("E50").Value = '  +2.23%  '

$ws.Range("D51").Value = '''20.02'
$ws.Range("E51").Value = '  -5.39%  '
